# Update LR-pair (NATMI) edge statistics with recomputed TPM-based values.
# Only the cells whose underlying ligand/receptor TPM changed are touched;
# everything else (headers, cluster labels, counts in columns A-D/K/L) is left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.003680666666666667
$ws.Range("H2").Value = 0.011042
$ws.Range("I2").Value = 0.005828946138996241
$ws.Range("J2").Value = 0.005828946138996241
$ws.Range("M2").Value = 5.273684
$ws.Range("N2").Value = 15.821052
$ws.Range("O2").Value = 0.0510821201937383
$ws.Range("P2").Value = 0.0510821201937383
$ws.Range("Q2").Value = 0.01941067290933333
$ws.Range("R2").Value = 0.174696056184
$ws.Range("S2").Value = 0.0002977549272750328
$ws.Range("T2").Value = 0.0002977549272750328

# Row 3
$ws.Range("G3").Value = 0.003680666666666667
$ws.Range("H3").Value = 0.011042
$ws.Range("I3").Value = 0.005828946138996241
$ws.Range("J3").Value = 0.005828946138996241
$ws.Range("O3").Value = 0.5598845502029881
$ws.Range("P3").Value = 0.5598845502029881
$ws.Range("Q3").Value = 0.2127502897248889
$ws.Range("R3").Value = 1.914752607524
$ws.Range("S3").Value = 0.003263536887189355
$ws.Range("T3").Value = 0.003263536887189354

# Row 4
$ws.Range("G4").Value = 0.003680666666666667
$ws.Range("H4").Value = 0.011042
$ws.Range("I4").Value = 0.005828946138996241
$ws.Range("J4").Value = 0.005828946138996241
$ws.Range("M4").Value = 32.95839133333334
$ws.Range("N4").Value = 98.87517400000002
$ws.Range("O4").Value = 0.3192425840231603
$ws.Range("P4").Value = 0.3192425840231604
$ws.Range("Q4").Value = 0.1213088523675556
$ws.Range("R4").Value = 1.091779671308
$ws.Range("S4").Value = 0.001860847827544983
$ws.Range("T4").Value = 0.001860847827544983

# Row 5
$ws.Range("G5").Value = 0.003680666666666667
$ws.Range("H5").Value = 0.011042
$ws.Range("I5").Value = 0.005828946138996241
$ws.Range("J5").Value = 0.005828946138996241
$ws.Range("M5").Value = 7.205150000000001
$ws.Range("N5").Value = 21.61545
$ws.Range("O5").Value = 0.06979074558011317
$ws.Range("P5").Value = 0.06979074558011318
$ws.Range("Q5").Value = 0.02651975543333334
$ws.Range("R5").Value = 0.2386777989
$ws.Range("S5").Value = 0.0004068064969868696
$ws.Range("T5").Value = 0.0004068064969868697

# Row 6
$ws.Range("I6").Value = 0.8498346916787334
$ws.Range("J6").Value = 0.8498346916787333
$ws.Range("M6").Value = 5.273684
$ws.Range("N6").Value = 15.821052
$ws.Range("O6").Value = 0.0510821201937383
$ws.Range("P6").Value = 0.0510821201937383
$ws.Range("Q6").Value = 2.8299906765
$ws.Range("R6").Value = 25.4699160885
$ws.Range("S6").Value = 0.04341135786514159
$ws.Range("T6").Value = 0.04341135786514158

# Row 7
$ws.Range("I7").Value = 0.8498346916787334
$ws.Range("J7").Value = 0.8498346916787333
$ws.Range("O7").Value = 0.5598845502029881
$ws.Range("P7").Value = 0.5598845502029881
$ws.Range("S7").Value = 0.4758093140974427
$ws.Range("T7").Value = 0.4758093140974426

# Row 8
$ws.Range("I8").Value = 0.8498346916787334
$ws.Range("J8").Value = 0.8498346916787333
$ws.Range("M8").Value = 32.95839133333334
$ws.Range("N8").Value = 98.87517400000002
$ws.Range("O8").Value = 0.3192425840231603
$ws.Range("P8").Value = 0.3192425840231604
$ws.Range("Q8").Value = 17.68629674925
$ws.Range("R8").Value = 159.17667074325
$ws.Range("S8").Value = 0.2713034229640446
$ws.Range("T8").Value = 0.2713034229640446

# Row 9
$ws.Range("I9").Value = 0.8498346916787334
$ws.Range("J9").Value = 0.8498346916787333
$ws.Range("M9").Value = 7.205150000000001
$ws.Range("N9").Value = 21.61545
$ws.Range("O9").Value = 0.06979074558011317
$ws.Range("P9").Value = 0.06979074558011318
$ws.Range("Q9").Value = 3.866463618750001
$ws.Range("R9").Value = 34.79817256875
$ws.Range("S9").Value = 0.0593105967521044
$ws.Range("T9").Value = 0.0593105967521044

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.09114066666666666
$ws.Range("H10").Value = 0.273422
$ws.Range("I10").Value = 0.1443363621822704
$ws.Range("J10").Value = 0.1443363621822704
$ws.Range("M10").Value = 5.273684
$ws.Range("N10").Value = 15.821052
$ws.Range("O10").Value = 0.0510821201937383
$ws.Range("P10").Value = 0.0510821201937383
$ws.Range("Q10").Value = 0.4806470755493333
$ws.Range("R10").Value = 4.325823679944
$ws.Range("S10").Value = 0.007373007401321681
$ws.Range("T10").Value = 0.007373007401321682

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.09114066666666666
$ws.Range("H11").Value = 0.273422
$ws.Range("I11").Value = 0.1443363621822704
$ws.Range("J11").Value = 0.1443363621822704
$ws.Range("O11").Value = 0.5598845502029881
$ws.Range("P11").Value = 0.5598845502029881
$ws.Range("Q11").Value = 5.268122597098222
$ws.Range("R11").Value = 47.413103373884
$ws.Range("S11").Value = 0.08081169921835606
$ws.Range("T11").Value = 0.08081169921835607

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.09114066666666666
$ws.Range("H12").Value = 0.273422
$ws.Range("I12").Value = 0.1443363621822704
$ws.Range("J12").Value = 0.1443363621822704
$ws.Range("M12").Value = 32.95839133333334
$ws.Range("N12").Value = 98.87517400000002
$ws.Range("O12").Value = 0.3192425840231603
$ws.Range("P12").Value = 0.3192425840231604
$ws.Range("Q12").Value = 3.003849758380889
$ws.Range("R12").Value = 27.034647825428
$ws.Range("S12").Value = 0.04607831323157076
$ws.Range("T12").Value = 0.04607831323157078

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.09114066666666666
$ws.Range("H13").Value = 0.273422
$ws.Range("I13").Value = 0.1443363621822704
$ws.Range("J13").Value = 0.1443363621822704
$ws.Range("M13").Value = 7.205150000000001
$ws.Range("N13").Value = 21.61545
$ws.Range("O13").Value = 0.06979074558011317
$ws.Range("P13").Value = 0.06979074558011318
$ws.Range("Q13").Value = 0.6566821744333333
$ws.Range("R13").Value = 5.910139569900001
$ws.Range("S13").Value = 0.0100733423310219
$ws.Range("T13").Value = 0.01007334233102191
